$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "https://d1eh9yux7w8iql.cloudfront.net/product_images/413837_34622d71-7e56-4018-b00e-26298dea5124.jpg"
$ws.Range("B1").Value = "iphone 12 pro max negro/grafito"

$ws.Range("A2").Value = "https://d1eh9yux7w8iql.cloudfront.net/product_images/413837_81c85355-e6d5-47d7-b060-05c4d8968965.jpg"
$ws.Range("B2").Value = "iphone 12 pro max negro/grafito"

$ws.Range("A3").Value = "https://d1eh9yux7w8iql.cloudfront.net/product_images/413837_025f5c2e-e75f-4a13-b25e-d92c7dfd7e0b.jpg"
$ws.Range("B3").Value = "iphone 12 pro max negro/grafito"

$ws.Range("A4").Value = "https://d1eh9yux7w8iql.cloudfront.net/product_images/413837_fa07ef38-b2a1-44ae-ba4d-d5bf31ac39d1.jpg"
$ws.Range("B4").Value = "iphone 12 pro max negro/grafito"
